$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-26 06:22:07"
$wsZh.Range("G4").Value = "2016-02-26 06:22:54"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-26 06:22:20"
$wsDe.Range("G4").Value = "2016-02-26 06:23:16"
